# C5-PowerPoint.pptx — apply the captured edit:
#   the table on the "SOURCES OF FINANCE" slide had its quick-style changed
#   from {3EB2E00B-AE47-4924-A38C-C832775981A0} to {11D17BB7-FFAA-4453-A754-CFE3E4952FA2}.

$p = $ppt.ActivePresentation

$targetStyle = "{11D17BB7-FFAA-4453-A754-CFE3E4952FA2}"
$applied = $false

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($targetStyle)
            $applied = $true
        }
    }
}

if (-not $applied) {
    Write-Host "WARNING: no table shape found to restyle"
}
